$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.800.12"
$ws.Range("E2").Value = "  +2.34%  "

$ws.Range("D3").Value = "2.300.88"
$ws.Range("E3").Value = "  +1.00%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").Value = "'322.68"
$ws.Range("E5").Value = "  +1.74%  "

$ws.Range("D6").Value = "'104.85"
$ws.Range("E6").Value = "  +2.84%  "

$ws.Range("E8").Value = "  +0.09%  "

$ws.Range("D9").Value = "'0.609"
$ws.Range("E9").Value = "  +1.15%  "

$ws.Range("D10").Value = "'40.19"
$ws.Range("E10").Value = "  +3.54%  "

$ws.Range("D11").Value = "'0.0910"
$ws.Range("E11").Value = "  +0.70%  "

$ws.Range("E12").Value = "  +3.84%  "

$ws.Range("E13").Value = "  +0.99%  "

$ws.Range("D14").Value = "'0.973"
$ws.Range("E14").Value = "  +1.96%  "

$ws.Range("D15").Value = "'15.31"
$ws.Range("E15").Value = "  +0.57%  "

$ws.Range("D16").Value = "2.649.94"
$ws.Range("E16").Value = "  +1.11%  "

$ws.Range("D17").Value = "2.292.23"
$ws.Range("E17").Value = "  +0.05%  "

$ws.Range("D18").Value = "42.727.64"
$ws.Range("E18").Value = "  +2.27%  "

$ws.Range("D19").Value = "'7.53"
$ws.Range("E19").Value = "  +0.55%  "

$ws.Range("E20").Value = "  +0.92%  "

$ws.Range("D21").Value = "'13.43"
$ws.Range("E21").Value = "  +35.02%  "

$ws.Range("D22").Value = "'73.66"
$ws.Range("E22").Value = "  +0.53%  "

$ws.Range("D24").Value = "'271.61"
$ws.Range("E24").Value = "  -3.17%  "

$ws.Range("E25").Value = "  -0.56%  "

$ws.Range("E26").Value = "  -0.40%  "

$ws.Range("D27").Value = "'10.95"
$ws.Range("E27").Value = "  +1.85%  "

$ws.Range("D28").Value = "'2.33"
$ws.Range("E28").Value = "  +0.74%  "

$ws.Range("D29").Value = "'22.67"
$ws.Range("E29").Value = "  -1.13%  "

$ws.Range("D30").Value = "'38.39"
$ws.Range("E30").Value = "  +11.79%  "

$ws.Range("D31").Value = "'165.64"
$ws.Range("E31").Value = "  +1.73%  "

$ws.Range("D32").Value = "'6.17"
$ws.Range("E32").Value = "  +6.10%  "

$ws.Range("D33").Value = "'0.0884"
$ws.Range("E33").Value = "  +1.57%  "

$ws.Range("E34").Value = "  +0.20%  "

$ws.Range("D35").Value = "'2.54"
$ws.Range("E35").Value = "  -12.42%  "

$ws.Range("E36").Value = "  -0.23%  "

$ws.Range("E37").Value = "  +1.32%  "

$ws.Range("E38").Value = "  +3.05%  "

$ws.Range("D39").Value = "'3.75"
$ws.Range("E39").Value = "  +4.00%  "

$ws.Range("E40").Value = "  -5.32%  "

$ws.Range("E41").Value = "  +6.69%  "

$ws.Range("D42").Value = "'99.19"
$ws.Range("E42").Value = "  -3.37%  "

$ws.Range("D43").Value = "'70.46"
$ws.Range("E43").Value = "  +1.88%  "

$ws.Range("E44").Value = "  -0.11%  "

$ws.Range("E45").Value = "  +0.47%  "

$ws.Range("D46").Value = "'12.45"
$ws.Range("E46").Value = "  +5.22%  "

$ws.Range("D47").Value = "'82.65"
$ws.Range("E47").Value = "  +9.33%  "

$ws.Range("D48").Value = "'113.64"
$ws.Range("E48").Value = "  -1.41%  "

$ws.Range("E49").Value = "  -1.17%  "

$ws.Range("D50").Value = "'5.30"
$ws.Range("E50").Value = "  +0.70%  "

$ws.Range("D51").Value = "1.593.90"
$ws.Range("E51").Value = "  +4.07%  "
